$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.273.18'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.000.63'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.72'
$ws.Range("E5").Value = '  -0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.54'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.436'
$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.51'
$ws.Range("E9").Value = '  -0.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("E11").Value = '  +3.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.522.19'
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("E13").Value = '  +1.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.35'
$ws.Range("E14").Value = '  +2.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000163'
$ws.Range("E15").Value = '  +5.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.255.56'
$ws.Range("E16").Value = '  +1.12%  '

$ws.Range("E17").Value = '  +6.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.999.62'
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.93'
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '328.46'
$ws.Range("E21").Value = '  +0.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.495'
$ws.Range("E23").Value = '  +3.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.42'
$ws.Range("E24").Value = '  +2.99%  '

$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0912'
$ws.Range("E27").Value = '  -0.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.74'
$ws.Range("E28").Value = '  +2.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  +5.52%  '

$ws.Range("E30").Value = '  +1.75%  '

$ws.Range("E31").Value = '  -5.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.56'
$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.69'
$ws.Range("E33").Value = '  +3.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.28'
$ws.Range("E34").Value = '  -1.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.85'
$ws.Range("E35").Value = '  +3.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.27'
$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0679'
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.28'
$ws.Range("E38").Value = '  +2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.036.52'
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.21'
$ws.Range("E40").Value = '  +1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.83'
$ws.Range("E42").Value = '  +5.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.291.19'
$ws.Range("E43").Value = '  +1.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.650'
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.983'
$ws.Range("E46").Value = '  -1.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.00'
$ws.Range("E47").Value = '  +3.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0238'
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.31'
$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("E50").Value = '  -7.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0891'
$ws.Range("E51").Value = '  +1.84%  '
